# "cambie los nombres de los excel, empleados ahora solo funciona con 1 excel"
# -> consolidate the sheet down to a single header row + a single data row,
#    give the header row real column titles instead of placeholder numbers,
#    and refresh the remaining data row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old sample rows (3-6); only one data row remains afterwards.
$ws.Rows("3:6").Delete()

# Row 1: turn the numeric placeholder headers into proper text labels.
$ws.Range("B1").Value = "Nombre"
$ws.Range("C1").Value = "Telefono"
$ws.Range("D1").Value = "Recurrencia"

# Row 2: updated record values (A stays numeric id, B stays the name
# "santiago", C becomes a numeric phone, D becomes a numeric recurrence).
$ws.Range("A2").Value = 1088
$ws.Range("B2").Value = "santiago"
$ws.Range("C2").Value = 3333
$ws.Range("D2").Value = 1
